$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1) Remove the block of bulleted "crossover" notes, from
#    "Once crossover node had been selected..." through
#    "Based on mutation rate which was relatively low"
# -----------------------------------------------------------------------
$count = $d.Paragraphs.Count
$startIdx = -1
$endIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "Once crossover node had been selected*") {
        $startIdx = $i
    }
    if ($t -like "Based on mutation rate which was relatively low*") {
        $endIdx = $i
        break
    }
}
if ($startIdx -gt 0 -and $endIdx -gt 0) {
    $startPar = $d.Paragraphs.Item($startIdx)
    $endPar = $d.Paragraphs.Item($endIdx)
    $rng = $d.Range($startPar.Range.Start, $endPar.Range.End)
    $rng.Delete()
}

# -----------------------------------------------------------------------
# 2) Remove the block of "prefix/infix parser" notes, from
#    "Now that the children trees were made..." through
#    "...then return the infix expression"
# -----------------------------------------------------------------------
$count = $d.Paragraphs.Count
$startIdx = -1
$endIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "Now that the children trees were made*") {
        $startIdx = $i
    }
    if ($t -like "It checks to see if the genetic operators are in the list*") {
        $endIdx = $i
        break
    }
}
if ($startIdx -gt 0 -and $endIdx -gt 0) {
    $startPar = $d.Paragraphs.Item($startIdx)
    $endPar = $d.Paragraphs.Item($endIdx)
    $rng = $d.Range($startPar.Range.Start, $endPar.Range.End)
    $rng.Delete()
}

# -----------------------------------------------------------------------
# 3) Move the "_GoBack" bookmark out of the final text paragraph and into
#    the (already existing) empty paragraph that follows it, so the
#    trailing space run stays with the sentence and the bookmark ends up
#    alone in its own paragraph.
# -----------------------------------------------------------------------
$count = $d.Paragraphs.Count
$pFinal = $d.Paragraphs.Item($count)

if ($d.Bookmarks.Exists("_GoBack")) {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
    $d.Bookmarks.Add("_GoBack", $pFinal.Range)
}
